$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.427.53"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.223.75"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.83"
$ws.Range("E5").Value = "  -8.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.52"
$ws.Range("E6").Value = "  +7.70%  "
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.52"
$ws.Range("E10").Value = "  -8.29%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.10"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.61"
$ws.Range("E13").Value = "  -8.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +10.58%  "
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("E16").Value = "  -5.61%  "
$ws.Range("D17").Value = "2.558.71"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "2.217.35"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "42.340.94"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("E21").Value = "  -4.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.72"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +13.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.68"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.96"
$ws.Range("E26").Value = "  -7.49%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -7.14%  "
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.60"
$ws.Range("E30").Value = "  -10.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.30"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.08"
$ws.Range("E32").Value = "  -8.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.84"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0878"
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0372"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("E40").Value = "  -4.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.22"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.39"
$ws.Range("E42").Value = "  -6.49%  "
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.24"
$ws.Range("E45").Value = "  -11.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.31"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.28"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -2.14%  "
